# Fix classification report sort order.
#
# The "Rating" rows were listed in plain alphabetical order (A, AA, AAA, B, BB, BBB,
# C, CC, CCC, D). The correct credit-rating order is AAA, AA, A, BBB, BB, B, CCC, CC, C, D.
# This swaps the data between row 2 <-> row 4 (A <-> AAA), row 5 <-> row 7 (B <-> BBB),
# and row 8 <-> row 10 (C <-> CCC). Rows 3 (AA), 6 (BB), 9 (CC) and 11 (D) already sit in
# the right place and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (as it should read after the sort-order fix), keyed by row number.
$target = @{
    2  = @("AAA", "0.8846", "0.9583", "0.9200", "24")
    4  = @("A",   "0.9604", "0.9327", "0.9463", "208")
    5  = @("BBB", "0.9669", "0.9669", "0.9669", "363")
    7  = @("B",   "0.9359", "0.9481", "0.9419", "154")
    8  = @("CCC", "0.8621", "0.9615", "0.9091", "26")
    10 = @("C",   "1.0000", "1.0000", "1.0000", "4")
}

foreach ($row in $target.Keys) {
    $values = $target[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $text = $values[$i]

        if ($col -eq 1) {
            # Rating column holds plain text (e.g. "AAA") - never looks numeric.
            $ws.Cells.Item($row, $col).Value2 = $text
        } else {
            # Precision/Recall/F1/Support columns contain numeric-looking text
            # (e.g. "0.8846", "24") that must stay stored as text, matching the
            # original inline-string cells. A leading apostrophe keeps Excel
            # from re-typing the cell as a number.
            $ws.Cells.Item($row, $col).Formula = "'" + $text
        }
    }
}
